$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append three new author names to column A, continuing after the
# existing last row (A13), consolidating entries from Wernsdorf's
# poetae minores. Values are written in this order so the shared
# string table picks up the new unique strings in the same sequence
# as the source workbook.
$ws.Range("A15").Value = "Vomanus?"
$ws.Range("A14").Value = "Ofilius?"
$ws.Range("A16").Value = "Speratus?"

# Update the active selection to reflect the next empty row, matching
# how Excel leaves the cursor after data entry.
$ws.Range("A17").Select()
